# Generate Report for Handback
# Applies:
#  1. Status text change: "Ready for handoff" -> "Handback transform failed"
#     (shows up on Overview sheet B3/C3 and on the zh-cn / de-de sheets C3)
#  2. Adds an "Error Detail" (column K) value on row 3 of the zh-cn sheet
#  3. Adds an "Error Detail" (column K) value on row 3 of the de-de sheet

$wb = $excel.ActiveWorkbook

$zhCnErrorDetail = "Handback file name: riseod0b.vju is different with handoff file name: 159678a5-0e23-466b-bc4d-f1d5710e2463.a04a12eb9a4fc7b752bdd66cd73add307ee8570b.zh-cn."
$deDeErrorDetail = "Handback file name: riseod0b.vju is different with handoff file name: 159678a5-0e23-466b-bc4d-f1d5710e2463.a04a12eb9a4fc7b752bdd66cd73add307ee8570b.de-de."

# 1. Update the "Ready for handoff" status text to "Handback transform failed"
#    wherever it currently appears across all worksheets.
foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    $nrows = $used.Rows.Count
    $ncols = $used.Columns.Count
    for ($r = 1; $r -le $nrows; $r++) {
        for ($c = 1; $c -le $ncols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            if ($cell.Text -eq "Ready for handoff") {
                $cell.Value = "Handback transform failed"
            }
        }
    }
}

# 2. zh-cn sheet: set Error Detail (column K) on row 3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("K3").Value = $zhCnErrorDetail

# 3. de-de sheet: set Error Detail (column K) on row 3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K3").Value = $deDeErrorDetail
